# Reviewed data for reference/consensus sequences
# Fill in nearest_upstream_orf (E) / nearest_downstream_orf (F) values
# for the rows that previously had blank entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows where upstream/downstream ORF are both "NK" (not known)
$nkRows = @(3, 4, 5, 8, 9, 39, 45)
foreach ($r in $nkRows) {
    $ws.Cells.Item($r, 5).Value = "NK"
    $ws.Cells.Item($r, 6).Value = "NK"
}

# Row 7 has specific gene/locus values
$ws.Cells.Item(7, 5).Value = "KLF8"
$ws.Cells.Item(7, 6).Value = "ENSACUG00000005807"
